$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns D and E keep their literal text representation
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.360.26'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '1.865.03'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('D4').Value = '1.020'
$ws.Range('E4').Value = '  -0.56%  '
$ws.Range('D5').Value = '316.72'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').Value = '1.019'
$ws.Range('E6').Value = '  +1.19%  '
$ws.Range('D7').Value = '0.5098'
$ws.Range('E7').Value = '  -0.81%  '
$ws.Range('D8').Value = '0.3954'
$ws.Range('E8').Value = '  +1.08%  '
$ws.Range('D9').Value = '0.08335'
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').Value = '42.00'
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('D11').Value = '1.107'
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('D12').Value = '6.226'
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('E13').Value = '  -1.12%  '
$ws.Range('D14').Value = '1.851.34'
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').Value = '1.020'
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').Value = '7.189'
$ws.Range('E16').Value = '  -1.43%  '
$ws.Range('D17').Value = '0.00001105'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = '90.70'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').Value = '0.06742'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('D20').Value = '1.020'
$ws.Range('E20').Value = '  +1.70%  '
$ws.Range('D21').Value = '17.64'
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('D22').Value = '5.948'
$ws.Range('E22').Value = '  -1.45%  '
$ws.Range('D23').Value = '28.412.16'
$ws.Range('E23').Value = '  +0.54%  '
$ws.Range('D24').Value = '11.14'
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').Value = '2.284'
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('D26').Value = '161.74'
$ws.Range('E26').Value = '  +1.34%  '
$ws.Range('D27').Value = '2.042.56'
$ws.Range('E27').Value = '  -1.58%  '
$ws.Range('D28').Value = '20.66'
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('D29').Value = '2.360'
$ws.Range('E29').Value = '  -4.57%  '
$ws.Range('D30').Value = '127.17'
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').Value = '0.1046'
$ws.Range('E31').Value = '  -1.08%  '
$ws.Range('D32').Value = '1.031'
$ws.Range('E32').Value = '  -0.91%  '
$ws.Range('D33').Value = '5.764'
$ws.Range('E33').Value = '  -1.14%  '
$ws.Range('D34').Value = '3.641'
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('D35').Value = '0.02420'
$ws.Range('E35').Value = '  -1.22%  '
$ws.Range('D36').Value = '0.06466'
$ws.Range('E36').Value = '  -1.83%  '
$ws.Range('D37').Value = '0.2179'
$ws.Range('E37').Value = '  -1.32%  '
$ws.Range('D38').Value = '8.845'
$ws.Range('E38').Value = '  -8.54%  '
$ws.Range('D39').Value = '1.269'
$ws.Range('E39').Value = '  +2.73%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.6405'
$ws.Range('E40').Value = '  -1.51%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '1.178'
$ws.Range('E41').Value = '  -1.88%  '
$ws.Range('D42').Value = '5.005'
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('D44').Value = '0.6001'
$ws.Range('E44').Value = '  -1.90%  '
$ws.Range('D45').Value = '12.96'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('E46').Value = '  +1.39%  '
$ws.Range('D47').Value = '1.219'
$ws.Range('E47').Value = '  -5.15%  '
$ws.Range('D48').Value = '1.987'
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('D49').Value = '121.79'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').Value = '1.201'
$ws.Range('E50').Value = '  -3.06%  '
$ws.Range('D51').Value = '0.06837'
$ws.Range('E51').Value = '  -1.75%  '
